# "rapport continue de stage" — split a few runs the way Word's live
# grammar/spell-checker does after a manual edit, and insert ", ou les deux"
# into the "Faire une autre page (...)" bullet.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "Accueil de l'équipe (briefing, compte rendu des taches à faire…)"
# -> split the single run into 3 runs around the word "briefing":
#    "Accueil de l'équipe (" | "briefing" | ", compte rendu des taches à faire…)"
# ---------------------------------------------------------------------
$rBriefing = $d.Content
$rBriefing.Find.Execute("briefing", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
# Touching Font forces Word to re-serialize this sub-range as its own run,
# matching the run-split produced by the original edit.
$rBriefing.Font.Bold = 1
$rBriefing.Font.Bold = 0

# ---------------------------------------------------------------------
# Edit 2: "Faire une autre page (que ce soit formulaire pour proposition
# de formation ou page de présentation)" gains ", ou les deux" right
# before the closing parenthesis, split into 3 runs:
#    "...page de présentation" | ", ou les deux" | ")"
# ---------------------------------------------------------------------
$rPage = $d.Content
$rPage.Find.Execute("page de présentation)", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
# Exclude the trailing ")" from the matched range so the new text lands
# right before it.
$rPage.End = $rPage.End - 1
$beforeLen = $rPage.End - $rPage.Start
$rPage.InsertAfter(", ou les deux")
# $rPage now spans "...page de présentation, ou les deux"; carve out a
# range over just the freshly inserted text and force it into its own run.
$insertedStart = $rPage.Start + $beforeLen
$rInserted = $d.Range($insertedStart, $rPage.End)
$rInserted.Font.Bold = 1
$rInserted.Font.Bold = 0
# ...and do the same for the closing parenthesis that now follows it.
$rCloseParen = $d.Range($rPage.End, $rPage.End + 1)
$rCloseParen.Font.Bold = 1
$rCloseParen.Font.Bold = 0

# ---------------------------------------------------------------------
# Edit 3: " concernant l'adresse mail réceptrice (création d'un nouveau
# compte google pour ça)" -> split the single run into 3 runs around the
# word "mail":
#    " concernant l'adresse " | "mail" | " réceptrice (...)"
# "mail" occurs more than once in the document, so first narrow down to
# the unique containing phrase, then search for "mail" within that range.
# ---------------------------------------------------------------------
$rAddrOuter = $d.Content
$rAddrOuter.Find.Execute("concernant l’adresse mail réceptrice", $true, `
    $false, $false, $false, $false, $true, 1, $false, "", 0)
$rMail = $rAddrOuter.Duplicate
$rMail.Find.Execute("mail", $true, $false, $false, $false, $false, $true, `
    1, $false, "", 0)
$rMail.Font.Bold = 1
$rMail.Font.Bold = 0
